$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Activate()

# Update the "multiplier cutoff" parameter (About!B12). This value feeds
# every dependent supply-curve formula on this sheet (B20:B121) as well as
# the mirrored formulas / chart caches on the CSC-CSCCCMvSoECBtY sheet, so
# changing it here cascades a full recalculation through the workbook.
$ws.Range("B12").Value = 0.15

# Reflect the author's selection state after editing the parameter.
$ws.Range("B12:B16").Select() | Out-Null
